$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 154, pushing existing rows 154-170 down to 156-172
$ws.Rows.Item(154).Resize(2).Insert()

# Row 154: new weekly price record
$ws.Cells.Item(154, 1).Value = 5
$ws.Cells.Item(154, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(154, 3).Value = "Maule"
$ws.Cells.Item(154, 4).Value = 45212
$ws.Cells.Item(154, 5).Value = 7
$ws.Cells.Item(154, 6).Value = "Fruta"
$ws.Cells.Item(154, 7).Value = 100107
$ws.Cells.Item(154, 8).Value = "Otros"
$ws.Cells.Item(154, 9).Value = 100107002
$ws.Cells.Item(154, 10).Value = "Chirimoya"
$ws.Cells.Item(154, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(154, 12).Value = "Primera"
$ws.Cells.Item(154, 13).Value = 180
$ws.Cells.Item(154, 14).Value = 20000
$ws.Cells.Item(154, 15).Value = 20000
$ws.Cells.Item(154, 16).Value = 20000
$ws.Cells.Item(154, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(154, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(154, 19).Value = 2000
$ws.Cells.Item(154, 20).Value = 10

# Row 155: new weekly price record
$ws.Cells.Item(155, 1).Value = 5
$ws.Cells.Item(155, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(155, 3).Value = "Maule"
$ws.Cells.Item(155, 4).Value = 45212
$ws.Cells.Item(155, 5).Value = 7
$ws.Cells.Item(155, 6).Value = "Fruta"
$ws.Cells.Item(155, 7).Value = 100107
$ws.Cells.Item(155, 8).Value = "Otros"
$ws.Cells.Item(155, 9).Value = 100107002
$ws.Cells.Item(155, 10).Value = "Chirimoya"
$ws.Cells.Item(155, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(155, 12).Value = "Segunda"
$ws.Cells.Item(155, 13).Value = 150
$ws.Cells.Item(155, 14).Value = 18000
$ws.Cells.Item(155, 15).Value = 18000
$ws.Cells.Item(155, 16).Value = 18000
$ws.Cells.Item(155, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(155, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(155, 19).Value = 1800
$ws.Cells.Item(155, 20).Value = 10
